# cv121192a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sub-header placeholder text in B2 ("unnamed: 1_level_1") is corrected
# to read "total", and the two blank spacer rows that separated the table
# sections ("situação do domicílio" and "grandes regiões e unidades da
# federação") are removed so the region/state data rows sit directly under
# their labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-labelled sub-header in row 2.
$ws.Range("B2").Value = "total"

# Remove the two blank separator rows (originally rows 5 and 8); deleting
# the lower one first keeps the row numbers of the other deletion stable.
$ws.Rows("8:8").EntireRow.Delete()
$ws.Rows("5:5").EntireRow.Delete()
